# Update the "想去人数" (want-to-go count) column F for rows 2-10 on both
# the "展览" sheet and the "全部类型" sheet (they mirror the same data).
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 53
    3  = 2109
    4  = 1604
    5  = 318
    6  = 1030
    7  = 489
    9  = 5668
    10 = 79
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
